# Updates cryptos list prices (D) and 1h volume change (E) columns
# per the scraped data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.406.87"
$ws.Range("E2").Value = "  +0.01%  "

$ws.Range("D3").Value = "1.572.44"
$ws.Range("E3").Value = "  +0.08%  "

$ws.Range("E5").Value = "  +0.22%  "

$ws.Range("D6").Value = "'290.50"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.25%  "

$ws.Range("E7").Value = "  +3.10%  "

$ws.Range("D8").Value = "'49.87"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.98%  "

$ws.Range("E9").Value = "  +1.26%  "

$ws.Range("E10").Value = "  -0.97%  "

$ws.Range("E11").Value = "  +0.85%  "

$ws.Range("E12").Value = "  +0.27%  "

$ws.Range("D13").Value = "'21.26"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.43%  "

$ws.Range("D14").Value = "'6.014"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.77%  "

$ws.Range("D15").Value = "'6.935"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.88%  "

$ws.Range("D16").Value = "1.567.70"
$ws.Range("E16").Value = "  +0.01%  "

$ws.Range("D17").Value = "'0.00001134"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.46%  "

$ws.Range("D18").Value = "'89.91"
$ws.Range("D18").ClearFormats()

$ws.Range("D19").Value = "'0.06758"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.28%  "

$ws.Range("E20").Value = "  +0.15%  "

$ws.Range("D21").Value = "'16.82"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.20%  "

$ws.Range("D22").Value = "'6.216"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.86%  "

$ws.Range("E23").Value = "  +0.21%  "

$ws.Range("D24").Value = "22.398.03"
$ws.Range("E24").Value = "  -0.08%  "

$ws.Range("D25").Value = "'2.419"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.10%  "

$ws.Range("D26").Value = "'2.719"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -9.41%  "

$ws.Range("E27").Value = "  +1.90%  "

$ws.Range("D28").Value = "'146.76"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.71%  "

$ws.Range("D29").Value = "'5.030"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.76%  "

$ws.Range("D30").Value = "'126.32"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.85%  "

$ws.Range("D31").Value = "1.746.71"
$ws.Range("E31").Value = "  +0.10%  "

$ws.Range("D32").Value = "'6.187"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.86%  "

$ws.Range("D33").Value = "'2.009"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.85%  "

$ws.Range("D34").Value = "'0.9938"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -5.22%  "

$ws.Range("D35").Value = "'10.02"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.32%  "

$ws.Range("E36").Value = "  +1.97%  "

$ws.Range("D37").Value = "'0.02549"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.38%  "

$ws.Range("D38").Value = "'0.2315"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.26%  "

$ws.Range("D39").Value = "'0.06572"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.63%  "

$ws.Range("D40").Value = "'1.328"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +6.08%  "

$ws.Range("D41").Value = "'5.477"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.77%  "

$ws.Range("D42").Value = "'0.6455"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.92%  "

$ws.Range("E43").Value = "  -2.78%  "

$ws.Range("D44").Value = "'14.17"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.30%  "

$ws.Range("D46").Value = "'3.797"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.61%  "

$ws.Range("D47").Value = "'0.6016"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.23%  "

$ws.Range("D48").Value = "'1.303"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +7.73%  "

$ws.Range("E49").Value = "  -2.52%  "

$ws.Range("D50").Value = "'125.36"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.37%  "

$ws.Range("D51").Value = "'0.07329"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.54%  "
